# Add slides for Multicore
# -------------------------------------------------------------
# Slide 31 ("Hardware Prefetcher Design Space") - body placeholder:
#   - Remove the leading "Predictors " run from the
#     "Predictors regular patterns (...)" bullet and capitalize
#     "regular" -> "Regular".
#   - Drop the leading "Predicted " word from the
#     "Predicted correlated patterns (...)" bullet and capitalize
#     "correlated" -> "Correlated".
# -------------------------------------------------------------

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(31)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Bullet: "Predictors regular patterns (x, x+8, x+16, ...)" ---
# The bullet is split across three runs: "Predictors ", "regular "
# (Wingdings bullet-glyph run) and "patterns (...)".  Re-capitalize
# the Wingdings-styled run first (its own text is unique so it is
# easy to find), then delete the now-stale "Predictors " run that
# precedes it.

$full = $tr.Text
$regularStart = $full.IndexOf("regular ") + 1
$regularRun = $tr.Characters($regularStart, 8)
$regularRun.Text = "Regular "

$full = $tr.Text
$predictorsStart = $full.IndexOf("Predictors ") + 1
$predictorsRun = $tr.Characters($predictorsStart, 11)
$predictorsRun.Text = ""

# --- Bullet: "Predicted correlated patterns (A...B->C, ...)" ---
# This bullet is a single run; replace the "Predicted " prefix and
# capitalize "correlated" -> "Correlated" in one shot.

$full = $tr.Text
$bulletStart = $full.IndexOf("Predicted correlated patterns") + 1
$oldBulletText = "Predicted correlated patterns (A" + [char]0x2026 + "B->C, B..C->J, A..C->K, " + [char]0x2026 + ")"
$newBulletText = "Correlated patterns (A" + [char]0x2026 + "B->C, B..C->J, A..C->K, " + [char]0x2026 + ")"
$bulletRun = $tr.Characters($bulletStart, $oldBulletText.Length)
$bulletRun.Text = $newBulletText

Write-Host "Slide 31 body text now:"
Write-Host $tr.Text
